$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for the Price/Volume columns so numeric-looking
# strings (e.g. "30.563.72", "1.000") are preserved as text, matching
# the source data which stores these as inline strings, not numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.563.72"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.914.80"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "245.10"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").Value = "0.4851"
$ws.Range("E7").Value = "  +2.19%  "
$ws.Range("D8").Value = "0.2890"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "0.06810"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "111.30"
$ws.Range("E10").Value = "  +5.85%  "
$ws.Range("E11").Value = "  +5.42%  "
$ws.Range("D12").Value = "1.918.77"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "0.07583"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "5.392"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "0.6719"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "294.40"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "30.559.02"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "13.03"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").Value = "0.000007602"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "5.524"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Value = "2.162.18"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "6.432"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "9.476"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "166.17"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").Value = "20.28"
$ws.Range("E27").Value = "  -4.05%  "
$ws.Range("D28").Value = "2.087"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "0.1065"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").Value = "1.439"
$ws.Range("E30").Value = "  +2.83%  "
$ws.Range("D31").Value = "4.125"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").Value = "0.7343"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "1.143"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "0.9994"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02035"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("D39").Value = "2.682"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "2.021"
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("D41").Value = "109.28"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").Value = "0.4440"
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("D43").Value = "0.8675"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "5.831"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "69.45"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").Value = "7.214"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "48.48"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "9.264"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "0.1227"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "0.2510"
$ws.Range("E51").Value = "  +0.52%  "
